$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 9.172748
$ws.Range("N2").Value = 27.518244
$ws.Range("O2").Value = 0.01445826353606064
$ws.Range("P2").Value = 0.01445826353606064
$ws.Range("Q2").Value = 0.519562792216
$ws.Range("R2").Value = 4.676065129944
$ws.Range("S2").Value = 0.01445826353606064
$ws.Range("T2").Value = 0.01445826353606064

# Row 3
$ws.Range("O3").Value = 0.2254554169720557
$ws.Range("P3").Value = 0.2254554169720557
$ws.Range("S3").Value = 0.2254554169720557
$ws.Range("T3").Value = 0.2254554169720557

# Row 4
$ws.Range("M4").Value = 169.2367096666667
$ws.Range("N4").Value = 507.7101290000001
$ws.Range("O4").Value = 0.2667541884216647
$ws.Range("P4").Value = 0.2667541884216647
$ws.Range("Q4").Value = 9.585905708939334
$ws.Range("R4").Value = 86.27315138045401
$ws.Range("S4").Value = 0.2667541884216647
$ws.Range("T4").Value = 0.2667541884216647

# Row 5
$ws.Range("M5").Value = 16.15031566666667
$ws.Range("N5").Value = 48.450947
$ws.Range("O5").Value = 0.02545644119943506
$ws.Range("P5").Value = 0.02545644119943505
$ws.Range("Q5").Value = 0.9147861799913334
$ws.Range("R5").Value = 8.233075619921999
$ws.Range("S5").Value = 0.02545644119943506
$ws.Range("T5").Value = 0.02545644119943505

# Row 6
$ws.Range("M6").Value = 54.744643
$ws.Range("N6").Value = 164.233929
$ws.Range("O6").Value = 0.08628956945961638
$ws.Range("P6").Value = 0.08628956945961638
$ws.Range("Q6").Value = 3.100846068806
$ws.Range("R6").Value = 27.907614619254
$ws.Range("S6").Value = 0.08628956945961638
$ws.Range("T6").Value = 0.08628956945961638

# Row 7
$ws.Range("M7").Value = 242.0894676666667
$ws.Range("N7").Value = 726.268403
$ws.Range("O7").Value = 0.3815861204111676
$ws.Range("P7").Value = 0.3815861204111676
$ws.Range("Q7").Value = 13.71243162757533
$ws.Range("R7").Value = 123.411884648178
$ws.Range("S7").Value = 0.3815861204111676
$ws.Range("T7").Value = 0.3815861204111676
